$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update user records (password / account details refresh)
$ws.Range("A2").Value = "newuser"
$ws.Range("B2").Value = "Mohammad"
$ws.Range("C2").Value = "Irfan"
$ws.Range("D2").Value = "irfan22@gmail.com"

$ws.Range("A3").Value = "newuser1"
$ws.Range("B3").Value = "Mohammad"
$ws.Range("C3").Value = "Farhan"
$ws.Range("D3").Value = "farhan@gmil.com"

# Move selection to H4
$ws.Range("H4").Select()
